# Weekly data refresh: a new daily price record is inserted as row 52
# ("Fruta / hortaliza, semanal"), pushing the existing rows 52-79 down
# to 53-80 (the sheet grows from A1:R79 to A1:R80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..79 down to 53..80, leaving a blank row 52 to fill in.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new market record.
$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 45001
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = 100112043
$ws.Cells.Item(52, 7).Value = "Pepino dulce"
$ws.Cells.Item(52, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 20000
$ws.Cells.Item(52, 12).Value = 20000
$ws.Cells.Item(52, 13).Value = 20000
$ws.Cells.Item(52, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(52, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 16).Value = 1111
$ws.Cells.Item(52, 17).Value = 18
$ws.Cells.Item(52, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date number format as the
# rest of column D (the Insert above already carries the format down,
# this is just a safety net).
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
